# Applies:
#  - "survey" sheet, cell M7: "instanceName" -> "instance_name"
#  - "survey" sheet, cell N7: TRUE -> FALSE
#  - "survey" sheet view: drop the G1 topLeftCell scroll and move the
#    selection from J8 to E7

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("survey")

$ws.Range("M7").Value = "instance_name"
$ws.Range("N7").Value = $false

$ws.Activate()
$ws.Range("E7").Select()
